$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1404.1538
$ws.Range("J80").Value = 1352.4
$ws.Range("L80").Value = 4057.2
$ws.Range("N80").Value = -6053.200000000001

$ws.Range("H83").Value = 1404.1538
$ws.Range("J83").Value = 1352.4
$ws.Range("L83").Value = 12171.6
$ws.Range("N83").Value = -22155.6

$ws.Range("H132").Value = 7273.85
$ws.Range("I132").Value = 3607.9395
$ws.Range("K132").Value = 10823.8185
$ws.Range("M132").Value = -8293.818499999999

$ws.Range("H135").Value = 1039.4791
$ws.Range("I135").Value = 680.81396
$ws.Range("K135").Value = 6127.325639999999
$ws.Range("M135").Value = -3592.325639999999

$ws.Range("H137").Value = 5683.0884
$ws.Range("I137").Value = 10168.643
$ws.Range("K137").Value = 30505.929
$ws.Range("M137").Value = -27955.929

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3930.94
$ws.Range("I32").Value = 3449.2708
$ws.Range("K32").Value = 3449.2708
$ws.Range("M32").Value = -3162.2708

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").ClearContents()
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = 0

$ws.Range("H39").Value = 16
$ws.Range("I39").Value = 16
$ws.Range("K39").Value = 16
$ws.Range("M39").Value = 504

$ws.Range("H61").Value = 3525.3225
$ws.Range("I61").Value = 3260.2144
$ws.Range("J61").Value = 5999.6665
$ws.Range("K61").Value = 3260.2144
$ws.Range("L61").Value = 5999.6665
$ws.Range("M61").Value = -3048.2144
$ws.Range("N61").Value = -6423.6665

$ws.Range("H74").Value = 2278.3635
$ws.Range("J74").Value = 3749.25
$ws.Range("L74").Value = 3749.25
$ws.Range("N74").Value = -5497.25

$ws.Range("H77").Value = 2278.3635
$ws.Range("J77").Value = 3749.25
$ws.Range("L77").Value = 18746.25
$ws.Range("N77").Value = -27482.25

$ws.Range("H132").Value = 1287.65
$ws.Range("I132").Value = 1162.1714
$ws.Range("J132").Value = 2166
$ws.Range("K132").Value = 3486.5142
$ws.Range("L132").Value = 6498
$ws.Range("M132").Value = -956.5141999999996
$ws.Range("N132").Value = -11558

$ws.Range("H136").Value = 3525.3225
$ws.Range("I136").Value = 3260.2144
$ws.Range("J136").Value = 5999.6665
$ws.Range("K136").Value = 9780.643199999999
$ws.Range("L136").Value = 17998.9995
$ws.Range("M136").Value = -7230.643199999999
$ws.Range("N136").Value = -23098.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3068.9546
$ws.Range("I107").Value = 2053.7058
$ws.Range("J107").Value = 6520.8
$ws.Range("K107").Value = 2053.7058
$ws.Range("L107").Value = 6520.8
$ws.Range("M107").Value = -133.7058000000002
$ws.Range("N107").Value = -10360.8

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").ClearContents()
$ws.Range("N119").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2790.0293
$ws.Range("I16").Value = 2241.6667
$ws.Range("J16").Value = 3675.8462
$ws.Range("K16").Value = 2241.6667
$ws.Range("L16").Value = 3675.8462
$ws.Range("M16").Value = -1954.6667
$ws.Range("N16").Value = -4249.8462

$ws.Range("H31").Value = 1776.8474
$ws.Range("I31").Value = 1440.6957
$ws.Range("K31").Value = 1440.6957
$ws.Range("M31").Value = -1145.6957

$ws.Range("H34").Value = 1776.8474
$ws.Range("I34").Value = 1440.6957
$ws.Range("K34").Value = 1440.6957
$ws.Range("M34").Value = -1238.6957

$ws.Range("H58").Value = 1374.7241
$ws.Range("I58").Value = 1283.9546
$ws.Range("J58").Value = 1660
$ws.Range("K58").Value = 1283.9546
$ws.Range("L58").Value = 1660
$ws.Range("M58").Value = -1080.9546
$ws.Range("N58").Value = -2066

$ws.Range("H113").Value = 2790.0293
$ws.Range("I113").Value = 2241.6667
$ws.Range("J113").Value = 3675.8462
$ws.Range("K113").Value = 2241.6667
$ws.Range("L113").Value = 3675.8462
$ws.Range("M113").Value = -71.66670000000022
$ws.Range("N113").Value = -8015.8462

$ws.Range("H136").Value = 1374.7241
$ws.Range("I136").Value = 1283.9546
$ws.Range("J136").Value = 1660
$ws.Range("K136").Value = 3851.8638
$ws.Range("L136").Value = 4980
$ws.Range("M136").Value = -1301.8638
$ws.Range("N136").Value = -10080

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 11365.444
$ws.Range("I10").Value = 18995.375
$ws.Range("K10").Value = 56986.125
$ws.Range("M10").Value = -56847.125

$ws.Range("H109").Value = 1629.8889
$ws.Range("I109").Value = 1666.9412
$ws.Range("K109").Value = 5000.8236
$ws.Range("M109").Value = -3960.8236

$ws.Range("H128").Value = 199933.47
$ws.Range("I128").Value = 199933.47
$ws.Range("K128").Value = 599800.41
$ws.Range("M128").Value = -594820.41

$ws.Range("H129").Value = 1157.1333
$ws.Range("J129").Value = 1503
$ws.Range("L129").Value = 4509
$ws.Range("N129").Value = -14509

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 12692.954
$ws.Range("I99").Value = 4022.6365
$ws.Range("K99").Value = 4022.6365
$ws.Range("M99").Value = -1776.6365

$ws.Range("H122").Value = 2285.4644
$ws.Range("J122").Value = 2999.75
$ws.Range("L122").Value = 8999.25
$ws.Range("N122").Value = -13899.25

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4412.3105
$ws.Range("I61").Value = 5190.2915
$ws.Range("K61").Value = 5190.2915
$ws.Range("M61").Value = -4988.2915

$ws.Range("H99").Value = 64999
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H113").Value = 4412.3105
$ws.Range("I113").Value = 5190.2915
$ws.Range("K113").Value = 5190.2915
$ws.Range("M113").Value = -3020.2915

$ws.Range("H132").Value = 14144.8
$ws.Range("I132").Value = 15843.907
$ws.Range("J132").Value = 3707.4285
$ws.Range("K132").Value = 47531.721
$ws.Range("L132").Value = 11122.2855
$ws.Range("M132").Value = -45001.721
$ws.Range("N132").Value = -16182.2855

$ws.Range("H136").Value = 5145183.5
$ws.Range("I136").Value = 6001085
$ws.Range("K136").Value = 18003255
$ws.Range("M136").Value = -18000705

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 27777.4
$ws.Range("I74").Value = 29705
$ws.Range("J74").Value = 26951.285
$ws.Range("K74").Value = 29705
$ws.Range("L74").Value = 26951.285
$ws.Range("M74").Value = -28769
$ws.Range("N74").Value = -28823.285

$ws.Range("H75").Value = 27498.5
$ws.Range("I75").Value = 24998
$ws.Range("J75").Value = 29999
$ws.Range("K75").Value = 24998
$ws.Range("L75").Value = 29999
$ws.Range("M75").Value = -24062
$ws.Range("N75").Value = -31871

$ws.Range("H77").Value = 27777.4
$ws.Range("I77").Value = 29705
$ws.Range("J77").Value = 26951.285
$ws.Range("K77").Value = 89115
$ws.Range("L77").Value = 80853.855
$ws.Range("M77").Value = -84435
$ws.Range("N77").Value = -90213.855

$ws.Range("H78").Value = 27498.5
$ws.Range("I78").Value = 24998
$ws.Range("J78").Value = 29999
$ws.Range("K78").Value = 74994
$ws.Range("L78").Value = 89997
$ws.Range("M78").Value = -70314
$ws.Range("N78").Value = -99357

$ws.Range("H126").Value = 8937134
$ws.Range("I126").Value = 9623952
$ws.Range("J126").Value = 8500
$ws.Range("K126").Value = 28871856
$ws.Range("L126").Value = 25500
$ws.Range("M126").Value = -28869386
$ws.Range("N126").Value = -30440

$ws.Range("H132").Value = 2165.6558
$ws.Range("I132").Value = 1623.7455
$ws.Range("J132").Value = 7133.1665
$ws.Range("K132").Value = 4871.2365
$ws.Range("L132").Value = 21399.4995
$ws.Range("M132").Value = -2341.2365
$ws.Range("N132").Value = -26459.4995
